# Generate Report for Handoff
# Updates status text, timestamps, and narrows the "datetime" columns.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Status text: "Handed back: in sync with en-US" -> "Ready for handoff" ---
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$zhcn.Range("C2").Value = "Ready for handoff"
$dede.Range("C2").Value = "Ready for handoff"

# --- Timestamps ---
# Overview!G2 (Latest HO Xliff Generate Date) and de-de!H2 (Latest Handoff Datetime)
$overview.Range("G2").Value = "2016-08-17 08:58:28"
$dede.Range("H2").Value = "2016-08-17 08:58:28"

# zh-cn!H2 (Latest Handoff Datetime)
$zhcn.Range("H2").Value = "2016-08-17 08:58:23"

# --- Column width changes (characters) ---
$overview.Columns.Item(5).ColumnWidth = 17.2159881591797
$overview.Columns.Item(6).ColumnWidth = 17.2159881591797
$zhcn.Columns.Item(3).ColumnWidth = 17.2159881591797
$dede.Columns.Item(3).ColumnWidth = 17.2159881591797
